# Update to skill section.
#
# 1. Drop the stray _GoBack bookmark that previously sat after the
#    "Experience in all mainstream ..." bullet.
# 2. The "Graphics programming and 3D mathematics (demos available)" bullet
#    becomes "Proficiency in developing a wide range of interactive
#    technology and simulations" (built from four separate runs).
# 3. The "4+ years developing interactive technology and simulations" bullet
#    becomes "Graphics programming and 3D mathematics (demos available)"
#    (built from two runs with the _GoBack bookmark re-inserted between them).

$d = $word.ActiveDocument

function Find-ParagraphByText($doc, $targetText) {
    foreach ($p in $doc.Paragraphs) {
        $t = $p.Range.Text
        $t = $t.TrimEnd("`r", "`n", "`a")
        if ($t -eq $targetText) {
            return $p
        }
    }
    return $null
}

# --- Step 1: remove the old _GoBack bookmark -----------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$rPrFonts = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr>'

$xmlHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' + `
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>'
$xmlFooter = '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# --- Step 2: "Graphics programming and 3D mathematics (demos available)" -
#     -> four runs of "Proficiency in developing a wide range of
#        interactive technology and simulations"
$pGraphics = Find-ParagraphByText $d "Graphics programming and 3D mathematics (demos available)"
if ($pGraphics -ne $null) {
    $full = $pGraphics.Range
    $target = $d.Range($full.Start, $full.End - 1)

    $body = ('<w:r w:rsidRPr="0053511D">' + $rPrFonts + '<w:t xml:space="preserve">Proficiency in </w:t></w:r>') + `
            ('<w:r>' + $rPrFonts + '<w:t xml:space="preserve">developing </w:t></w:r>') + `
            ('<w:r>' + $rPrFonts + '<w:t xml:space="preserve">a wide range of </w:t></w:r>') + `
            ('<w:r>' + $rPrFonts + '<w:t>interactive technology and simulations</w:t></w:r>')

    $target.InsertXML($xmlHeader + $body + $xmlFooter)
}

# --- Step 3: "4+ years developing interactive technology and simulations" -
#     -> "Graphics program" + _GoBack bookmark + "ming and 3D mathematics
#        (demos available)"
$pYears = Find-ParagraphByText $d "4+ years developing interactive technology and simulations"
if ($pYears -ne $null) {
    $full = $pYears.Range
    $target = $d.Range($full.Start, $full.End - 1)

    $body = ('<w:r w:rsidRPr="0053511D">' + $rPrFonts + '<w:t>Graphics program</w:t></w:r>') + `
            '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' + `
            ('<w:r>' + $rPrFonts + '<w:t>ming and 3D mathematics (demos available)</w:t></w:r>')

    $target.InsertXML($xmlHeader + $body + $xmlFooter)
}

Write-Output "skills section updated"
